# Updated capital structure database
# Refreshes the south_korea_banks_regional sheet's financial metrics
# (rows 2-5) with newer source data, including two company rows whose
# names swap places (row 3 <-> row 5) and a few cells that no longer
# have values in the refreshed dataset (F3, AN2, AP2, AN3, AP3).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (South Korea / Banks (Regional), company #1)
$ws.Range("D2").Value = 0.0369
$ws.Range("E2").Value = -0.0316
$ws.Range("F2").Value = 0.031
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 717
$ws.Range("L2").Value = 0.1258490864093518
$ws.Range("M2").Value = 168.5334
$ws.Range("N2").Value = 0.05920307724733903
$ws.Range("O2").Value = 0.2350535564853556
$ws.Range("P2").Value = 162.5234
$ws.Range("Q2").Value = 0.05709186075104507
$ws.Range("R2").Value = 0.2266714086471408
$ws.Range("S2").Value = 6.010000000000005
$ws.Range("T2").Value = 0.03566058715957789
$ws.Range("U2").Value = 2683.1
$ws.Range("V2").Value = 0.9425299469561247
$ws.Range("W2").Value = 0.05992998968727759
$ws.Range("X2").Value = 0.2280416679541334
$ws.Range("Y2").Value = -0.1681116782668558
$ws.Range("Z2").Value = 0.160971599064227
$ws.Range("AB2").Value = 0.03867922227675859
$ws.Range("AC2").Value = -0.03867922227675859
$ws.Range("AD2").Value = 35036.4
$ws.Range("AE2").Value = 0
$ws.Range("AF2").Value = 35036.4
$ws.Range("AG2").Value = 32353.3
$ws.Range("AH2").Value = 0.9248556744300229
$ws.Range("AI2").Value = 0.7193845180922032
$ws.Range("AJ2").Value = 0.9191278409090909
$ws.Range("AK2").Value = 0.7030238895093893

# Row 3 (company renamed from DGB Financial Group to Jeju Bank)
$ws.Range("B3").Value = 'Jeju Bank (KOSE:A006220)'
$ws.Range("D3").Value = 0.0369
$ws.Range("E3").Value = -0.0316
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 19.6
$ws.Range("L3").Value = 0.1658206429780034
$ws.Range("M3").Value = 2.7606
$ws.Range("N3").Value = 0.02825588536335721
$ws.Range("O3").Value = 0.1408469387755102
$ws.Range("P3").Value = 2.7606
$ws.Range("Q3").Value = 0.02825588536335721
$ws.Range("R3").Value = 0.1408469387755102
$ws.Range("U3").Value = 55.6
$ws.Range("V3").Value = 0.5690890481064483
$ws.Range("W3").Value = 0.04770017035775129
$ws.Range("X3").Value = 0.1241071097666179
$ws.Range("Y3").Value = -0.07640693940886663
$ws.Range("Z3").Value = 0.1846586470863928
$ws.Range("AA3").Value = 0
$ws.Range("AB3").Value = 0.0385495859479186
$ws.Range("AC3").Value = -0.0385495859479186
$ws.Range("AD3").Value = 400.7
$ws.Range("AE3").Value = 0
$ws.Range("AF3").Value = 400.7
$ws.Range("AG3").Value = 345.1
$ws.Range("AH3").Value = 0.8039727126805779
$ws.Range("AI3").Value = 0.4774216609078994
$ws.Range("AJ3").Value = 0.7793586269196026
$ws.Range("AK3").Value = 0.440347071583514

# Row 4
$ws.Range("D4").Value = -0.012
$ws.Range("E4").Value = -0.129
$ws.Range("F4").Value = 0.031
$ws.Range("K4").Value = 412.6
$ws.Range("L4").Value = 0.1371629932515541
$ws.Range("M4").Value = 106.2496
$ws.Range("N4").Value = 0.06273965160909359
$ws.Range("O4").Value = 0.2575123606398449
$ws.Range("P4").Value = 100.2396
$ws.Range("Q4").Value = 0.05919078830823737
$ws.Range("R4").Value = 0.2429461948618516
$ws.Range("S4").Value = 6.010000000000005
$ws.Range("T4").Value = 0.05656491883263565
$ws.Range("U4").Value = 1117.6
$ws.Range("V4").Value = 0.6599350457632123
$ws.Range("W4").Value = 0.05992998968727759
$ws.Range("X4").Value = 0.2280416679541334
$ws.Range("Y4").Value = -0.1681116782668558
$ws.Range("Z4").Value = 0.1782915871454142
$ws.Range("AB4").Value = 0.03867922227675859
$ws.Range("AC4").Value = -0.03867922227675859
$ws.Range("AD4").Value = 15277.6
$ws.Range("AF4").Value = 15277.6
$ws.Range("AG4").Value = 14160
$ws.Range("AH4").Value = 0.9002127145559217
$ws.Range("AI4").Value = 0.6498755769189868
$ws.Range("AJ4").Value = 0.8931781625508562
$ws.Range("AK4").Value = 0.6323997695492365

# Row 5 (company renamed from Jeju Bank to DGB Financial Group)
$ws.Range("B5").Value = 'DGB Financial Group Co., Ltd. (KOSE:A139130)'
$ws.Range("D5").Value = 0.07820000000000001
$ws.Range("E5").Value = 0.0133
$ws.Range("K5").Value = 284.8
$ws.Range("L5").Value = 0.1107740178918709
$ws.Range("M5").Value = 59.5232
$ws.Range("N5").Value = 0.05639336807200378
$ws.Range("O5").Value = 0.209
$ws.Range("P5").Value = 59.5232
$ws.Range("Q5").Value = 0.05639336807200378
$ws.Range("R5").Value = 0.209
$ws.Range("U5").Value = 1509.9
$ws.Range("V5").Value = 1.430506868782567
$ws.Range("W5").Value = 0.07408756275851305
$ws.Range("X5").Value = 0.4249031693490773
$ws.Range("Y5").Value = -0.3508156065905642
$ws.Range("Z5").Value = 0.1437814923970852
$ws.Range("AB5").Value = 0.03874398866061474
$ws.Range("AC5").Value = -0.03874398866061474
$ws.Range("AD5").Value = 19358.1
$ws.Range("AF5").Value = 19358.1
$ws.Range("AG5").Value = 17848.2
$ws.Range("AH5").Value = 0.9482942744052985
$ws.Range("AI5").Value = 0.7948143129888525
$ws.Range("AJ5").Value = 0.9441643699381602
$ws.Range("AK5").Value = 0.7812532829078684

# --- Clear cells removed in the update ---
$ws.Range("AN2").ClearContents()
$ws.Range("AP2").ClearContents()
$ws.Range("AN3").ClearContents()
$ws.Range("AP3").ClearContents()
$ws.Range("F3").ClearContents()

Write-Output "Applied south_korea_banks_regional updates"